$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the containing path "sequence/run_4019_samples/" from each fastqFileName
# value in column F (rows 2-19), leaving just the bare filename.
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $current = $cell.Value2
    if ($current -ne $null) {
        $newValue = $current -replace '^sequence/run_4019_samples/', ''
        $cell.Value = $newValue
    }
}

# Update the selection to reflect the new active range (was C2:C19, now F2:F19)
$ws.Range("F2:F19").Select()
